$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update styling of B2/C2 (Liverpool / Real Madrid) so they match the
# "plain" (no font color) centered style already used by D2 ---
$ws.Range("D2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C2").PasteSpecial(-4122)

# --- Add new row 3: Game Id 2, Barcelona vs Paris, date 10/12/2014 ---

# A3 - Game Id, reuse style of A2 ("themed" style)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 2

# B3 - Home Team "Barcelona", reuse the original "themed" style (still on A2)
$ws.Range("A2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = "Barcelona"

# C3 - Away Team "Paris", reuse the "plain" style (still on D2)
$ws.Range("D2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "Paris"

# D3 - Date 10/12/2014 stored as serial 41983, "plain" style + built-in
# short-date number format (numFmtId 14)
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = 41983
$ws.Range("D3").NumberFormat = "mm-dd-yy"

$ws.Rows.Item(3).RowHeight = 14.25

$excel.CutCopyMode = 0
